$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D (Price) and E (Volume) columns are treated as plain text so values
# like "1.00" or "0.0685" keep their exact literal formatting instead of being
# auto-coerced into numbers by the COM layer.

# Row 2: update D2, E2
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = '42.254.22'
$ws.Range("E2").Value = '  -0.80%  '

# Row 3: update D3, E3
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = '2.271.03'
$ws.Range("E3").Value = '  -1.07%  '

# Row 4: update D4, E4
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.01%  '

# Row 5: update D5, E5
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = '299.92'
$ws.Range("E5").Value = '  -0.99%  '

# Row 6: update D6, E6
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = '95.96'
$ws.Range("E6").Value = '  -2.70%  '

# Row 7: update E7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.38%  '

# Row 8: update E8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.03%  '

# Row 9: update D9, E9
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = '0.491'
$ws.Range("E9").Value = '  -1.85%  '

# Row 10: update D10, E10
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = '33.21'
$ws.Range("E10").Value = '  -3.85%  '

# Row 11: update D11, E11
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0788'
$ws.Range("E11").Value = '  -0.19%  '

# Row 12: update D12, E12
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = '48.40'
$ws.Range("E12").Value = '  -6.46%  '

# Row 13: update D13, E13
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = '0.114'
$ws.Range("E13").Value = '  +0.97%  '

# Row 14: update D14, E14
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = '6.66'
$ws.Range("E14").Value = '  -1.24%  '

# Row 15: update D15, E15
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = '15.64'
$ws.Range("E15").Value = '  -0.57%  '

# Row 16: update D16, E16
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = '2.623.31'
$ws.Range("E16").Value = '  -1.08%  '

# Row 17: update D17, E17
$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = '2.281.88'
$ws.Range("E17").Value = '  -0.29%  '

# Row 18: update D18, E18
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = '0.785'
$ws.Range("E18").Value = '  -2.23%  '

# Row 19: update D19, E19
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = '42.184.63'
$ws.Range("E19").Value = '  -0.76%  '

# Row 20: update D20, E20
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = '11.72'
$ws.Range("E20").Value = '  +2.09%  '

# Row 21: update E21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.25%  '

# Row 22: update E22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.14%  '

# Row 23: update E23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.53%  '

# Row 24: update D24, E24
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = '235.02'
$ws.Range("E24").Value = '  -0.05%  '

# Row 25: update E25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.27%  '

# Row 26: update E26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.01%  '

# Row 27: update D27, E27
$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = '2.46'
$ws.Range("E27").Value = '  -2.30%  '

# Row 28: update D28, E28
$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D28").Value = '23.87'
$ws.Range("E28").Value = '  -4.49%  '

# Row 29: update E29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.13%  '

# Row 30: update D30, E30
$ws.Range("D30:E30").NumberFormat = "@"
$ws.Range("D30").Value = '168.00'
$ws.Range("E30").Value = '  +2.77%  '

# Row 31: update D31, E31
$ws.Range("D31:E31").NumberFormat = "@"
$ws.Range("D31").Value = '9.19'
$ws.Range("E31").Value = '  +0.13%  '

# Row 32: update D32, E32
$ws.Range("D32:E32").NumberFormat = "@"
$ws.Range("D32").Value = '33.70'
$ws.Range("E32").Value = '  -3.06%  '

# Row 33: update D33, E33
$ws.Range("D33:E33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.02%  '

# Row 34: update D34, E34
$ws.Range("D34:E34").NumberFormat = "@"
$ws.Range("D34").Value = '4.89'
$ws.Range("E34").Value = '  -2.49%  '

# Row 35: update D35, E35
$ws.Range("D35:E35").NumberFormat = "@"
$ws.Range("D35").Value = '4.55'
$ws.Range("E35").Value = '  -1.37%  '

# Row 36: update E36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.15%  '

# Row 37: update E37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.30%  '

# Row 38: update D38, E38
$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0685'
$ws.Range("E38").Value = '  -3.72%  '

# Row 39: update D39, E39
$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = '2.80'
$ws.Range("E39").Value = '  -2.83%  '

# Row 40: update D40, E40
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0986'
$ws.Range("E40").Value = '  -1.83%  '

# Row 41: update B41, C41, D41, E41
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = '1.72'
$ws.Range("E41").Value = '  -4.39%  '

# Row 42: update B42, C42, D42, E42
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").Value = '0.109'
$ws.Range("E42").Value = '  -2.49%  '

# Row 43: update E43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.17%  '

# Row 44: update D44, E44
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = '1.968.14'
$ws.Range("E44").Value = '  -0.71%  '

# Row 45: update E45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.06%  '

# Row 46: update D46, E46
$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = '17.48'
$ws.Range("E46").Value = '  -5.70%  '

# Row 47: update E47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -6.75%  '

# Row 48: update D48, E48
$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = '2.78'
$ws.Range("E48").Value = '  -4.37%  '

# Row 49: update D49, E49
$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = '2.494.98'
$ws.Range("E49").Value = '  -1.00%  '

# Row 50: update D50, E50
$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = '52.39'
$ws.Range("E50").Value = '  -5.49%  '

# Row 51: update E51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.51%  '
